$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the duplicate bold "Play Age of the Gods Furious 4 for
#    Free | Review" paragraph that was sitting near the end of the doc.
#    (Done first, before the Meta description paragraph is inserted at
#    the top, so there is only one paragraph with this exact text and
#    no ambiguity when locating it.)
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Age of the Gods Furious 4 for Free | Review" -and $p.Range.Bold) {
        $dupRange = $p.Range
        $dupRange.Expand(4)          # wdParagraph - include the paragraph mark
        $dupRange.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) Replace the old meta-description text (now italic, at the very end
#    of the document) with the DALLE image prompt, keeping formatting.
#    (Also done before the duplicate "Read our review..." text is
#    introduced by the new Meta description paragraph below.)
# ---------------------------------------------------------------------
$oldText = "Read our review of Age of the Gods Furious 4 and play this slot for free. Discover unique powers, progressive jackpots, and more exciting features."
$newText = 'Prompt for DALLE: Create a cartoon-style feature image for "Age of the Gods: Furious 4" online slot game. The image should feature a happy Maya warrior with glasses. Keep in mind the game''s theme of ancient Greek mythology and the four different deities: Prometheus, Apollo, Pandora, and Atlas. Use vibrant colors and bold lines to capture the excitement of the game''s features and special functions, such as free spins and progressive jackpots. Make sure to incorporate the game''s logo into the image.'

$searchRange = $d.Content
$found = $searchRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Setting .Text directly (vs. Find.Execute's Replace) keeps the run's
    # formatting (italic) and avoids smart-quote autocorrection.
    $searchRange.Text = $newText
}

# ---------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph at the top of the document.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)              # wdCollapseEnd
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"           # match body-text paragraphs (no pStyle)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Age of the Gods Furious 4 and play this slot for free. Discover unique powers, progressive jackpots, and more exciting features.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
